$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 15, shifting the existing row 15
# (summary_statistics_assembly) down to row 16.
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with the checksum field definition.
$ws.Range("A15").Value = "cheksum"
$ws.Range("B15").Value = "md5 checksum of the summary stats file"
$ws.Range("C15").Value = $true
$ws.Range("D15").Value = $false
$ws.Range("E15").Value = "string"
$ws.Range("J15").Value = "md5 sum"

# Update the selection to match the saved workbook view state.
$ws.Range("A15").Select()
